# Auto-generated Excel COM-interop script
# Updates market-board derived price/profit columns (H-N) on several
# crafting-leve worksheets to reflect refreshed Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 126.07143
$ws.Range("I5").Value = 130.07692
$ws.Range("K5").Value = 130.07692
$ws.Range("M5").Value = -15.07692
$ws.Range("H6").Value = 791.8570999999999
$ws.Range("I6").Value = 775.8461
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 2327.5383
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -2215.5383
$ws.Range("N6").Value = -3224
$ws.Range("H17").Value = 427488.4
$ws.Range("J17").Value = 427488.4
$ws.Range("L17").Value = 1282465.2
$ws.Range("N17").Value = -1282801.2
$ws.Range("H19").Value = 1247.8572
$ws.Range("I19").Value = 1305.2222
$ws.Range("J19").Value = 1144.6
$ws.Range("K19").Value = 1305.2222
$ws.Range("L19").Value = 1144.6
$ws.Range("M19").Value = -1130.2222
$ws.Range("N19").Value = -1494.6
$ws.Range("H39").Value = 1331.3
$ws.Range("I39").Value = 1289.375
$ws.Range("K39").Value = 3868.125
$ws.Range("M39").Value = -3572.125
$ws.Range("H55").Value = 80
$ws.Range("I55").Value = 73.333336
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 73.333336
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = 140.666664
$ws.Range("N55").Value = -528
$ws.Range("H70").Value = 10141
$ws.Range("I70").Value = 2729.3333
$ws.Range("J70").Value = 11851.385
$ws.Range("K70").Value = 8187.999899999999
$ws.Range("L70").Value = 35554.155
$ws.Range("M70").Value = -7917.999899999999
$ws.Range("N70").Value = -36094.155
$ws.Range("H73").Value = 10141
$ws.Range("I73").Value = 2729.3333
$ws.Range("J73").Value = 11851.385
$ws.Range("K73").Value = 8187.999899999999
$ws.Range("L73").Value = 35554.155
$ws.Range("M73").Value = -7251.999899999999
$ws.Range("N73").Value = -37426.155
$ws.Range("H74").Value = 11985.3125
$ws.Range("I74").Value = 11985.3125
$ws.Range("K74").Value = 11985.3125
$ws.Range("M74").Value = -11049.3125
$ws.Range("H77").Value = 11985.3125
$ws.Range("I77").Value = 11985.3125
$ws.Range("K77").Value = 59926.5625
$ws.Range("M77").Value = -55246.5625
$ws.Range("H106").Value = 2219.8333
$ws.Range("I106").Value = 1904.7778
$ws.Range("J106").Value = 3165
$ws.Range("K106").Value = 1904.7778
$ws.Range("L106").Value = 3165
$ws.Range("M106").Value = -1273.7778
$ws.Range("N106").Value = -4427
$ws.Range("H132").Value = 2713.0527
$ws.Range("I132").Value = 2713.0527
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8139.158100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5609.158100000001
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 17857656
$ws.Range("I135").Value = 18519014
$ws.Range("K135").Value = 166671126
$ws.Range("M135").Value = -166668591
$ws.Range("H141").Value = 754.3
$ws.Range("I141").Value = 745.8276
$ws.Range("K141").Value = 2237.4828
$ws.Range("M141").Value = 2942.5172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1246.6492
$ws.Range("I32").Value = 1232.5741
$ws.Range("K32").Value = 1232.5741
$ws.Range("M32").Value = -945.5741

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2293.6875
$ws.Range("I16").Value = 1389.7
$ws.Range("K16").Value = 1389.7
$ws.Range("M16").Value = -1102.7
$ws.Range("H58").Value = 20838900
$ws.Range("I58").Value = 25006014
$ws.Range("K58").Value = 25006014
$ws.Range("M58").Value = -25005811
$ws.Range("H113").Value = 2293.6875
$ws.Range("I113").Value = 1389.7
$ws.Range("K113").Value = 1389.7
$ws.Range("M113").Value = 780.3
$ws.Range("H122").Value = 3859.7778
$ws.Range("I122").Value = 3859.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11579.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9129.3334
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 20834982
$ws.Range("I132").Value = 25642288
$ws.Range("K132").Value = 76926864
$ws.Range("M132").Value = -76924334
$ws.Range("H136").Value = 20838900
$ws.Range("I136").Value = 25006014
$ws.Range("K136").Value = 75018042
$ws.Range("M136").Value = -75015492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1101
$ws.Range("H113").Value = 167449.67
$ws.Range("I113").Value = 500950.5
$ws.Range("J113").Value = 699.25
$ws.Range("K113").Value = 1502851.5
$ws.Range("L113").Value = 2097.75
$ws.Range("M113").Value = -1500681.5
$ws.Range("N113").Value = -6437.75
$ws.Range("H132").Value = 2200
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19800
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -24860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4313749.5
$ws.Range("I132").Value = 5955394
$ws.Range("J132").Value = 4431.625
$ws.Range("K132").Value = 17866182
$ws.Range("L132").Value = 13294.875
$ws.Range("M132").Value = -17863652
$ws.Range("N132").Value = -18354.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8934932
$ws.Range("I132").Value = 10006692
$ws.Range("J132").Value = 3598
$ws.Range("K132").Value = 30020076
$ws.Range("L132").Value = 10794
$ws.Range("M132").Value = -30017546
$ws.Range("N132").Value = -15854
$ws.Range("H136").Value = 2024.2693
$ws.Range("I136").Value = 1331.6471
$ws.Range("K136").Value = 3994.9413
$ws.Range("M136").Value = -1444.9413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2139.5833
$ws.Range("I81").Value = 2152.0908
$ws.Range("J81").Value = 2002
$ws.Range("K81").Value = 4304.1816
$ws.Range("L81").Value = 4004
$ws.Range("M81").Value = -3243.1816
$ws.Range("N81").Value = -6126
$ws.Range("H84").Value = 2139.5833
$ws.Range("I84").Value = 2152.0908
$ws.Range("J84").Value = 2002
$ws.Range("K84").Value = 21520.908
$ws.Range("L84").Value = 20020
$ws.Range("M84").Value = -16216.908
$ws.Range("N84").Value = -30628
$ws.Range("H122").Value = 1244.8572
$ws.Range("I122").Value = 1244.8572
$ws.Range("K122").Value = 3734.5716
$ws.Range("M122").Value = -1284.5716
$ws.Range("H132").Value = 9617457
$ws.Range("I132").Value = 12822602
$ws.Range("K132").Value = 38467806
$ws.Range("M132").Value = -38465276
$ws.Range("H136").Value = 10418898
$ws.Range("I136").Value = 10871676
$ws.Range("K136").Value = 32615028
$ws.Range("M136").Value = -32612478
